# Nayem_meal.xlsx edit: Rakib +1000, Rakib Extra meal, Bazar 11, Nayem +11
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rakib had an extra meal on day 5 (F3): 1.5 -> 2.5
$ws.Range("F3").Value = 2.5

# Rakib deposited/was credited +1000 on day 6 of the cost/deposit table (G23): 0 -> 1000
$ws.Range("G23").Value = 1000

# Bazar (market) spend of 11 recorded on day 6 of the cost/deposit table (G29): 0 -> 11
$ws.Range("G29").Value = 11

# New "Nayem" entry added to the per-person Bazar table (row 42 = names, row 43 = amounts)
$ws.Range("F42").Value = "Nayem"
$ws.Range("F43").Value = 11

# Restore the view state: scrolled so row 4 is the top-left row, with F45 selected
$ws.Range("F45").Select()

$wb.Application.Calculate()
